$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two new blank rows before the current row 12 ("endOfTestData"),
# pushing it down to row 14 and growing the used range to A1:D14.
$ws.Rows("12:13").Insert()

# Update the existing "Rahul Jha" / "Java Trainer" test-data row.
$ws.Range("A11").Value = "Rahul Arora"
$ws.Range("B11").Value = "Selenium Master"

# Fill in the two newly inserted rows with additional user test data.
$ws.Range("A12").Value = "Deepak"
$ws.Range("B12").Value = "Ui Developer"

$ws.Range("A13").Value = "Anil"
$ws.Range("B13").Value = "Tech Lead"

# Match the author's final selection.
$ws.Range("B13").Select()
